{"js": "// Find the paragraph that ends the document body text (\"Un experimento con\n// ITI grande evaluamos la capacidad para retener el evento anterior.\") and\n// insert five new \"Normal\" paragraphs right after it:\n//   1. (empty paragraph)\n//   2. \"Que busca una ratas? \"\n//   3. \"Comer el m\u00e1ximo en menor tiempo\"\n//   4. \"Comer el m\u00e1ximo posible\"\n//   5. \"Comer en el menor tiempo posible\"\n// The two pre-existing trailing empty paragraphs stay after the new block.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst marker = \"Un experimento con ITI grande evaluamos la capacidad para retener el evento anterior.\";\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === marker) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not locate the anchor paragraph for the insertion.\");\n}\n\nconst newTexts = [\n  \"\",\n  \"Que busca una ratas? \",\n  \"Comer el m\u00e1ximo en menor tiempo\",\n  \"Comer el m\u00e1ximo posible\",\n  \"Comer en el menor tiempo posible\"\n];\n\nlet anchor = target;\nfor (const text of newTexts) {\n  anchor = anchor.insertParagraph(text, Word.InsertLocation.after);\n}\n\nawait context.sync();\n", "ps1": "# Locate the paragraph that ends with:\n#   \"Un experimento con ITI grande evaluamos la capacidad para retener el evento anterior.\"\n# and insert five new \"Normal\" paragraphs right after it:\n#   1. (empty paragraph)\n#   2. \"Que busca una ratas? \"\n#   3. \"Comer el m\u00e1ximo en menor tiempo\"\n#   4. \"Comer el m\u00e1ximo posible\"\n#   5. \"Comer en el menor tiempo posible\"\n# The two pre-existing trailing empty paragraphs remain after the new block.\n\n$d = $word.ActiveDocument\n\n$marker = \"Un experimento con ITI grande evaluamos la capacidad para retener el evento anterior.\"\n\n# Find the anchor text and remember where it starts.\n$findRange = $d.Content\n$findRange.Find.Execute($marker) | Out-Null\n$anchorStart = $findRange.Start\n\n# Resolve the paragraph index (1-based) whose range starts there.\n$count = $d.Paragraphs.Count\n$anchorIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Start -eq $anchorStart) {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -eq -1) {\n    throw \"Could not locate the anchor paragraph for the insertion.\"\n}\n\n$newTexts = @(\n    \"\",\n    \"Que busca una ratas? \",\n    \"Comer el m\u00e1ximo en menor tiempo\",\n    \"Comer el m\u00e1ximo posible\",\n    \"Comer en el menor tiempo posible\"\n)\n\n$curIndex = $anchorIndex\nforeach ($text in $newTexts) {\n    $para = $d.Paragraphs.Item($curIndex)\n    $para.Range.InsertParagraphAfter()\n    $curIndex = $curIndex + 1\n    if ($text -ne \"\") {\n        $newPara = $d.Paragraphs.Item($curIndex)\n        $newPara.Range.Text = $text\n    }\n}\n"}
